$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Total")
$ws.Activate()

$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 4
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 8

$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = 5
$ws.Range("H3").Value = 8
$ws.Range("I3").Value = 4

$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 4
$ws.Range("G4").Value = 4
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 4

$ws.Range("E5").Value = 4
$ws.Range("F5").Value = 4
$ws.Range("G5").Value = 4
$ws.Range("H5").Value = 4
$ws.Range("I5").Value = 4

$ws.Range("M2").Select()
$ws.Application.ActiveWindow.ScrollColumn = 2
